$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-08 Wednesday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-03-09 Thursday", 2) | Out-Null
$d.Content.Find.Execute("45×95=4275", $true, $true, $false, $false, $false, $true, 1, $false, "50×17=850", 2) | Out-Null
$d.Content.Find.Execute("58×50=2900", $true, $true, $false, $false, $false, $true, 1, $false, "59×77=4543", 2) | Out-Null
$d.Content.Find.Execute("65×31=2015", $true, $true, $false, $false, $false, $true, 1, $false, "97×29=2813", 2) | Out-Null
$d.Content.Find.Execute("24×77=1848", $true, $true, $false, $false, $false, $true, 1, $false, "16×93=1488", 2) | Out-Null
$d.Content.Find.Execute("98×52=5096", $true, $true, $false, $false, $false, $true, 1, $false, "50×65=3250", 2) | Out-Null
$d.Content.Find.Execute("37×33=1221", $true, $true, $false, $false, $false, $true, 1, $false, "80×68=5440", 2) | Out-Null
$d.Content.Find.Execute("67×31=2077", $true, $true, $false, $false, $false, $true, 1, $false, "98×23=2254", 2) | Out-Null
$d.Content.Find.Execute("11×47=517", $true, $true, $false, $false, $false, $true, 1, $false, "98×71=6958", 2) | Out-Null
$d.Content.Find.Execute("75×35=2625", $true, $true, $false, $false, $false, $true, 1, $false, "54×97=5238", 2) | Out-Null
$d.Content.Find.Execute("41×89=3649", $true, $true, $false, $false, $false, $true, 1, $false, "80×18=1440", 2) | Out-Null
$d.Content.Find.Execute("23×11=253", $true, $true, $false, $false, $false, $true, 1, $false, "58×24=1392", 2) | Out-Null
$d.Content.Find.Execute("63×64=4032", $true, $true, $false, $false, $false, $true, 1, $false, "56×15=840", 2) | Out-Null
$d.Content.Find.Execute("95×32=3040", $true, $true, $false, $false, $false, $true, 1, $false, "23×100=2300", 2) | Out-Null
$d.Content.Find.Execute("97×14=1358", $true, $true, $false, $false, $false, $true, 1, $false, "70×12=840", 2) | Out-Null
$d.Content.Find.Execute("40×88=3520", $true, $true, $false, $false, $false, $true, 1, $false, "54×58=3132", 2) | Out-Null
$d.Content.Find.Execute("72×11=792", $true, $true, $false, $false, $false, $true, 1, $false, "87×61=5307", 2) | Out-Null
$d.Content.Find.Execute("45×33=1485", $true, $true, $false, $false, $false, $true, 1, $false, "55×83=4565", 2) | Out-Null
$d.Content.Find.Execute("29×16=464", $true, $true, $false, $false, $false, $true, 1, $false, "99×12=1188", 2) | Out-Null
$d.Content.Find.Execute("50×23=1150", $true, $true, $false, $false, $false, $true, 1, $false, "13×57=741", 2) | Out-Null
$d.Content.Find.Execute("60×79=4740", $true, $true, $false, $false, $false, $true, 1, $false, "82×52=4264", 2) | Out-Null
$d.Content.Find.Execute("43×64=2752", $true, $true, $false, $false, $false, $true, 1, $false, "51×78=3978", 2) | Out-Null
$d.Content.Find.Execute("18×22=396", $true, $true, $false, $false, $false, $true, 1, $false, "66×86=5676", 2) | Out-Null
$d.Content.Find.Execute("51×59=3009", $true, $true, $false, $false, $false, $true, 1, $false, "21×59=1239", 2) | Out-Null
$d.Content.Find.Execute("93×69=6417", $true, $true, $false, $false, $false, $true, 1, $false, "92×31=2852", 2) | Out-Null
$d.Content.Find.Execute("100×70=7000", $true, $true, $false, $false, $false, $true, 1, $false, "89×74=6586", 2) | Out-Null
$d.Content.Find.Execute("33×95=3135", $true, $true, $false, $false, $false, $true, 1, $false, "26×56=1456", 2) | Out-Null
$d.Content.Find.Execute("11×21=231", $true, $true, $false, $false, $false, $true, 1, $false, "50×28=1400", 2) | Out-Null
$d.Content.Find.Execute("90×99=8910", $true, $true, $false, $false, $false, $true, 1, $false, "27×11=297", 2) | Out-Null
$d.Content.Find.Execute("69×54=3726", $true, $true, $false, $false, $false, $true, 1, $false, "21×31=651", 2) | Out-Null
$d.Content.Find.Execute("38×60=2280", $true, $true, $false, $false, $false, $true, 1, $false, "88×86=7568", 2) | Out-Null
$d.Content.Find.Execute("64×27=1728", $true, $true, $false, $false, $false, $true, 1, $false, "30×19=570", 2) | Out-Null
$d.Content.Find.Execute("75×23=1725", $true, $true, $false, $false, $false, $true, 1, $false, "89×27=2403", 2) | Out-Null
$d.Content.Find.Execute("54×56=3024", $true, $true, $false, $false, $false, $true, 1, $false, "98×11=1078", 2) | Out-Null
$d.Content.Find.Execute("85×92=7820", $true, $true, $false, $false, $false, $true, 1, $false, "32×23=736", 2) | Out-Null
$d.Content.Find.Execute("83×55=4565", $true, $true, $false, $false, $false, $true, 1, $false, "14×90=1260", 2) | Out-Null
$d.Content.Find.Execute("63×84=5292", $true, $true, $false, $false, $false, $true, 1, $false, "51×41=2091", 2) | Out-Null
$d.Content.Find.Execute("81×33=2673", $true, $true, $false, $false, $false, $true, 1, $false, "57×44=2508", 2) | Out-Null
$d.Content.Find.Execute("78×85=6630", $true, $true, $false, $false, $false, $true, 1, $false, "34×68=2312", 2) | Out-Null
$d.Content.Find.Execute("43×96=4128", $true, $true, $false, $false, $false, $true, 1, $false, "34×26=884", 2) | Out-Null
$d.Content.Find.Execute("86×56=4816", $true, $true, $false, $false, $false, $true, 1, $false, "53×86=4558", 2) | Out-Null
$d.Content.Find.Execute("25×62=1550", $true, $true, $false, $false, $false, $true, 1, $false, "84×33=2772", 2) | Out-Null
$d.Content.Find.Execute("58×82=4756", $true, $true, $false, $false, $false, $true, 1, $false, "64×45=2880", 2) | Out-Null
$d.Content.Find.Execute("76×17=1292", $true, $true, $false, $false, $false, $true, 1, $false, "88×92=8096", 2) | Out-Null
$d.Content.Find.Execute("81×72=5832", $true, $true, $false, $false, $false, $true, 1, $false, "56×14=784", 2) | Out-Null
$d.Content.Find.Execute("22×54=1188", $true, $true, $false, $false, $false, $true, 1, $false, "79×18=1422", 2) | Out-Null
$d.Content.Find.Execute("81×15=1215", $true, $true, $false, $false, $false, $true, 1, $false, "18×98=1764", 2) | Out-Null
$d.Content.Find.Execute("41×46=1886", $true, $true, $false, $false, $false, $true, 1, $false, "72×78=5616", 2) | Out-Null
$d.Content.Find.Execute("35×40=1400", $true, $true, $false, $false, $false, $true, 1, $false, "39×74=2886", 2) | Out-Null
$d.Content.Find.Execute("62×14=868", $true, $true, $false, $false, $false, $true, 1, $false, "43×59=2537", 2) | Out-Null
$d.Content.Find.Execute("97×80=7760", $true, $true, $false, $false, $false, $true, 1, $false, "48×41=1968", 2) | Out-Null
$d.Content.Find.Execute("96×18=1728", $true, $true, $false, $false, $false, $true, 1, $false, "82×22=1804", 2) | Out-Null
$d.Content.Find.Execute("27×30=810", $true, $true, $false, $false, $false, $true, 1, $false, "35×20=700", 2) | Out-Null
$d.Content.Find.Execute("57×71=4047", $true, $true, $false, $false, $false, $true, 1, $false, "36×32=1152", 2) | Out-Null
$d.Content.Find.Execute("43×17=731", $true, $true, $false, $false, $false, $true, 1, $false, "11×43=473", 2) | Out-Null
$d.Content.Find.Execute("70×82=5740", $true, $true, $false, $false, $false, $true, 1, $false, "29×81=2349", 2) | Out-Null
$d.Content.Find.Execute("94×98=9212", $true, $true, $false, $false, $false, $true, 1, $false, "40×41=1640", 2) | Out-Null
$d.Content.Find.Execute("22×19=418", $true, $true, $false, $false, $false, $true, 1, $false, "69×90=6210", 2) | Out-Null
$d.Content.Find.Execute("84×68=5712", $true, $true, $false, $false, $false, $true, 1, $false, "88×58=5104", 2) | Out-Null
$d.Content.Find.Execute("87×58=5046", $true, $true, $false, $false, $false, $true, 1, $false, "39×45=1755", 2) | Out-Null
$d.Content.Find.Execute("19×97=1843", $true, $true, $false, $false, $false, $true, 1, $false, "15×29=435", 2) | Out-Null
$d.Content.Find.Execute("19×88=1672", $true, $true, $false, $false, $false, $true, 1, $false, "92×97=8924", 2) | Out-Null
$d.Content.Find.Execute("16×79=1264", $true, $true, $false, $false, $false, $true, 1, $false, "35×19=665", 2) | Out-Null
$d.Content.Find.Execute("100×67=6700", $true, $true, $false, $false, $false, $true, 1, $false, "85×57=4845", 2) | Out-Null
$d.Content.Find.Execute("61×30=1830", $true, $true, $false, $false, $false, $true, 1, $false, "94×63=5922", 2) | Out-Null
$d.Content.Find.Execute("100×64=6400", $true, $true, $false, $false, $false, $true, 1, $false, "75×61=4575", 2) | Out-Null
$d.Content.Find.Execute("22×92=2024", $true, $true, $false, $false, $false, $true, 1, $false, "82×89=7298", 2) | Out-Null
$d.Content.Find.Execute("14×30=420", $true, $true, $false, $false, $false, $true, 1, $false, "51×52=2652", 2) | Out-Null
$d.Content.Find.Execute("15×87=1305", $true, $true, $false, $false, $false, $true, 1, $false, "98×84=8232", 2) | Out-Null
$d.Content.Find.Execute("90×47=4230", $true, $true, $false, $false, $false, $true, 1, $false, "42×15=630", 2) | Out-Null
$d.Content.Find.Execute("51×33=1683", $true, $true, $false, $false, $false, $true, 1, $false, "42×24=1008", 2) | Out-Null
$d.Content.Find.Execute("15×56=840", $true, $true, $false, $false, $false, $true, 1, $false, "40×14=560", 2) | Out-Null
$d.Content.Find.Execute("75×25=1875", $true, $true, $false, $false, $false, $true, 1, $false, "28×51=1428", 2) | Out-Null
$d.Content.Find.Execute("12×66=792", $true, $true, $false, $false, $false, $true, 1, $false, "89×77=6853", 2) | Out-Null
$d.Content.Find.Execute("98×67=6566", $true, $true, $false, $false, $false, $true, 1, $false, "40×83=3320", 2) | Out-Null
$d.Content.Find.Execute("28×96=2688", $true, $true, $false, $false, $false, $true, 1, $false, "46×62=2852", 2) | Out-Null
$d.Content.Find.Execute("81×39=3159", $true, $true, $false, $false, $false, $true, 1, $false, "90×22=1980", 2) | Out-Null
$d.Content.Find.Execute("51×100=5100", $true, $true, $false, $false, $false, $true, 1, $false, "30×62=1860", 2) | Out-Null
$d.Content.Find.Execute("46×66=3036", $true, $true, $false, $false, $false, $true, 1, $false, "98×54=5292", 2) | Out-Null
$d.Content.Find.Execute("67×85=5695", $true, $true, $false, $false, $false, $true, 1, $false, "22×99=2178", 2) | Out-Null
$d.Content.Find.Execute("52×87=4524", $true, $true, $false, $false, $false, $true, 1, $false, "58×70=4060", 2) | Out-Null
$d.Content.Find.Execute("47×83=3901", $true, $true, $false, $false, $false, $true, 1, $false, "81×37=2997", 2) | Out-Null
$d.Content.Find.Execute("86×88=7568", $true, $true, $false, $false, $false, $true, 1, $false, "86×46=3956", 2) | Out-Null
$d.Content.Find.Execute("40×35=1400", $true, $true, $false, $false, $false, $true, 1, $false, "45×71=3195", 2) | Out-Null
$d.Content.Find.Execute("65×45=2925", $true, $true, $false, $false, $false, $true, 1, $false, "31×65=2015", 2) | Out-Null
$d.Content.Find.Execute("28×99=2772", $true, $true, $false, $false, $false, $true, 1, $false, "64×41=2624", 2) | Out-Null
$d.Content.Find.Execute("25×44=1100", $true, $true, $false, $false, $false, $true, 1, $false, "14×15=210", 2) | Out-Null
$d.Content.Find.Execute("86×81=6966", $true, $true, $false, $false, $false, $true, 1, $false, "16×88=1408", 2) | Out-Null
$d.Content.Find.Execute("71×55=3905", $true, $true, $false, $false, $false, $true, 1, $false, "43×65=2795", 2) | Out-Null
$d.Content.Find.Execute("24×60=1440", $true, $true, $false, $false, $false, $true, 1, $false, "83×46=3818", 2) | Out-Null
$d.Content.Find.Execute("33×60=1980", $true, $true, $false, $false, $false, $true, 1, $false, "78×89=6942", 2) | Out-Null
$d.Content.Find.Execute("13×18=234", $true, $true, $false, $false, $false, $true, 1, $false, "95×64=6080", 2) | Out-Null
$d.Content.Find.Execute("16×54=864", $true, $true, $false, $false, $false, $true, 1, $false, "29×39=1131", 2) | Out-Null
$d.Content.Find.Execute("87×27=2349", $true, $true, $false, $false, $false, $true, 1, $false, "35×49=1715", 2) | Out-Null
$d.Content.Find.Execute("91×26=2366", $true, $true, $false, $false, $false, $true, 1, $false, "16×56=896", 2) | Out-Null
$d.Content.Find.Execute("10×31=310", $true, $true, $false, $false, $false, $true, 1, $false, "99×68=6732", 2) | Out-Null
$d.Content.Find.Execute("32×72=2304", $true, $true, $false, $false, $false, $true, 1, $false, "80×61=4880", 2) | Out-Null
$d.Content.Find.Execute("49×36=1764", $true, $true, $false, $false, $false, $true, 1, $false, "84×41=3444", 2) | Out-Null
$d.Content.Find.Execute("12×94=1128", $true, $true, $false, $false, $false, $true, 1, $false, "71×81=5751", 2) | Out-Null
$d.Content.Find.Execute("28×45=1260", $true, $true, $false, $false, $false, $true, 1, $false, "29×62=1798", 2) | Out-Null
$d.Content.Find.Execute("60×20=1200", $true, $true, $false, $false, $false, $true, 1, $false, "97×81=7857", 2) | Out-Null
